$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 12270
$ws1.Range("F12").Value = 320
$ws1.Range("F19").Value = 338
$ws1.Range("F22").Value = 280
$ws1.Range("F24").Value = 331
$ws1.Range("F28").Value = 263
$ws1.Range("F29").Value = 786
$ws1.Range("F30").Value = 1271

# Sheet "本地生活" (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 1908

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1908
$ws4.Range("F5").Value = 12270
$ws4.Range("F17").Value = 320
$ws4.Range("F24").Value = 338
$ws4.Range("F27").Value = 280
$ws4.Range("F32").Value = 331
$ws4.Range("F38").Value = 263
$ws4.Range("F40").Value = 786
$ws4.Range("F41").Value = 1271
